# Hide slide 2 (sldId 258, cId 246086491) and slide 10 (sldId 263, cId 312535264).
# Both slides keep their position in the deck; only their "show in slide show"
# flag changes, which PowerPoint persists as show="0" on <p:sld>.

$p = $ppt.ActivePresentation

$slide2 = $p.Slides.Item(2)
$slide2.SlideShowTransition.Hidden = 1

$slide10 = $p.Slides.Item(10)
$slide10.SlideShowTransition.Hidden = 1
